$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read the current "Tipo" header text (D1) and "single" value (D2) before moving them
$tipoHeader = $ws.Range("D1").Value2
$tipoValue = $ws.Range("D2").Value2

# Copy the header style from D1 into the new E1 cell (it will hold "Tipo")
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats

# Move "Tipo"/"single" one column to the right, into the new E column
$ws.Range("E1").Value = $tipoHeader
$ws.Range("E2").Value = $tipoValue

# Add the new "MAE" header in D1, reusing the same header style
$ws.Range("D1").Value = "MAE"

# Update the metric values and add the new MAE value in D2
$ws.Range("B2").Value = 0.2492319147989941
$ws.Range("C2").Value = 0.9951287559575412
$ws.Range("D2").Value = 0.3920540822269443
